# Update cryptos list (prices + 1h volume %) per Wed Mar 29 22:25:28 UTC 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.403.67"
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("D3").Value = "1.793.86"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'314.39"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.5424"
$ws.Range("E7").Value = "  +3.14%  "
$ws.Range("D8").Value = "'0.3837"
$ws.Range("E8").Value = "  +4.37%  "
$ws.Range("D9").Value = "'0.07574"
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").Value = "'42.47"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "'1.121"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'21.09"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "'6.185"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").Value = "'7.401"
$ws.Range("E15").Value = "  +7.03%  "
$ws.Range("D16").Value = "1.798.93"
$ws.Range("E16").Value = "  +1.78%  "
$ws.Range("D17").Value = "'91.74"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "'0.06458"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'17.32"
$ws.Range("E21").Value = "  +3.73%  "
$ws.Range("D22").Value = "'5.960"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "28.396.30"
$ws.Range("E23").Value = "  +4.12%  "
$ws.Range("D24").Value = "'11.34"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'2.119"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'159.53"
$ws.Range("E26").Value = "  +2.96%  "
$ws.Range("D27").Value = "'20.67"
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("D28").Value = "'2.397"
$ws.Range("E28").Value = "  +3.39%  "
$ws.Range("D29").Value = "2.002.04"
$ws.Range("E29").Value = "  +1.52%  "
$ws.Range("D30").Value = "'123.27"
$ws.Range("E30").Value = "  +2.18%  "
$ws.Range("D31").Value = "'1.117"
$ws.Range("E31").Value = "  +5.51%  "
$ws.Range("D32").Value = "'0.1021"
$ws.Range("E32").Value = "  +4.52%  "
$ws.Range("D33").Value = "'5.733"
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("D34").Value = "'3.701"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("D35").Value = "'0.2326"
$ws.Range("E35").Value = "  +15.28%  "
$ws.Range("D36").Value = "'0.06369"
$ws.Range("E36").Value = "  +6.86%  "
$ws.Range("D37").Value = "'0.02318"
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("D38").Value = "'5.140"
$ws.Range("E38").Value = "  +6.49%  "
$ws.Range("D39").Value = "'8.778"
$ws.Range("E39").Value = "  +8.97%  "
$ws.Range("D40").Value = "'11.61"
$ws.Range("E40").Value = "  +3.86%  "
$ws.Range("D41").Value = "'0.6396"
$ws.Range("E41").Value = "  +4.41%  "
$ws.Range("D42").Value = "'1.001"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'1.156"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("D44").Value = "'1.391"
$ws.Range("D45").Value = "'13.50"
$ws.Range("E45").Value = "  +2.91%  "
$ws.Range("D46").Value = "'0.5966"
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "'3.673"
$ws.Range("E47").Value = "  +1.41%  "
$ws.Range("D48").Value = "'125.98"
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("D49").Value = "'1.982"
$ws.Range("E49").Value = "  +5.79%  "
$ws.Range("D50").Value = "'1.149"
$ws.Range("E50").Value = "  +3.26%  "
$ws.Range("E51").Value = "  +3.03%  "
